$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "None"
$ws.Range("F3").Value = "None"
$ws.Range("F4").Value = "None"
$ws.Range("F5").Value = "None"

$ws.Range("F2").Select()
